$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.5206698179245
$ws.Range("B1").Value = 2.175567388534546
$ws.Range("C1").Value = 2.592531204223633
$ws.Range("D1").Value = 3.100714683532715
$ws.Range("E1").Value = 2.119130849838257
